$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.142057
$ws.Range("H2").Value = 3.426171
$ws.Range("I2").Value = 0.2487160836358648
$ws.Range("J2").Value = 0.2487160836358647
$ws.Range("M2").Value = 6.492744666666667
$ws.Range("N2").Value = 19.478234
$ws.Range("O2").Value = 0.284509526105254
$ws.Range("P2").Value = 0.2845095261052539
$ws.Range("Q2").Value = 7.415084495779334
$ws.Range("R2").Value = 66.735760462014
$ws.Range("S2").Value = 0.07076209508999459
$ws.Range("T2").Value = 0.07076209508999458

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.142057
$ws.Range("H3").Value = 3.426171
$ws.Range("I3").Value = 0.2487160836358648
$ws.Range("J3").Value = 0.2487160836358647
$ws.Range("O3").Value = 0.1833546924709238
$ws.Range("P3").Value = 0.1833546924709237
$ws.Range("Q3").Value = 4.778717099499001
$ws.Range("R3").Value = 43.00845389549101
$ws.Range("S3").Value = 0.04560326102762653
$ws.Range("T3").Value = 0.04560326102762653

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.142057
$ws.Range("H4").Value = 3.426171
$ws.Range("I4").Value = 0.2487160836358648
$ws.Range("J4").Value = 0.2487160836358647
$ws.Range("M4").Value = 1.79534
$ws.Range("N4").Value = 5.38602
$ws.Range("O4").Value = 0.07867109501782452
$ws.Range("P4").Value = 0.0786710950178245
$ws.Range("Q4").Value = 2.05038061438
$ws.Range("R4").Value = 18.45342552942
$ws.Range("S4").Value = 0.01956676664817831
$ws.Range("T4").Value = 0.0195667666481783

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.142057
$ws.Range("H5").Value = 3.426171
$ws.Range("I5").Value = 0.2487160836358648
$ws.Range("J5").Value = 0.2487160836358647
$ws.Range("M5").Value = 10.34844233333333
$ws.Range("N5").Value = 31.045327
$ws.Range("O5").Value = 0.4534646864059979
$ws.Range("P5").Value = 0.4534646864059979
$ws.Range("Q5").Value = 11.81851100587967
$ws.Range("R5").Value = 106.366599052917
$ws.Range("S5").Value = 0.1127839608700653
$ws.Range("T5").Value = 0.1127839608700653

# Row 6
$ws.Range("I6").Value = 0.2524673276986635
$ws.Range("J6").Value = 0.2524673276986635
$ws.Range("M6").Value = 6.492744666666667
$ws.Range("N6").Value = 19.478234
$ws.Range("O6").Value = 0.284509526105254
$ws.Range("P6").Value = 0.2845095261052539
$ws.Range("Q6").Value = 7.526922022662666
$ws.Range("R6").Value = 67.742298203964
$ws.Range("S6").Value = 0.07182935976060661
$ws.Range("T6").Value = 0.0718293597606066

# Row 7
$ws.Range("I7").Value = 0.2524673276986635
$ws.Range("J7").Value = 0.2524673276986635
$ws.Range("O7").Value = 0.1833546924709238
$ws.Range("P7").Value = 0.1833546924709237
$ws.Range("S7").Value = 0.04629106922914437
$ws.Range("T7").Value = 0.04629106922914437

# Row 8
$ws.Range("I8").Value = 0.2524673276986635
$ws.Range("J8").Value = 0.2524673276986635
$ws.Range("M8").Value = 1.79534
$ws.Range("N8").Value = 5.38602
$ws.Range("O8").Value = 0.07867109501782452
$ws.Range("P8").Value = 0.0786710950178245
$ws.Range("Q8").Value = 2.08130534588
$ws.Range("R8").Value = 18.73174811292
$ws.Range("S8").Value = 0.0198618811262778
$ws.Range("T8").Value = 0.01986188112627779

# Row 9
$ws.Range("I9").Value = 0.2524673276986635
$ws.Range("J9").Value = 0.2524673276986635
$ws.Range("M9").Value = 10.34844233333333
$ws.Range("N9").Value = 31.045327
$ws.Range("O9").Value = 0.4534646864059979
$ws.Range("P9").Value = 0.4534646864059979
$ws.Range("Q9").Value = 11.99676292507133
$ws.Range("R9").Value = 107.970866325642
$ws.Range("S9").Value = 0.1144850175826347
$ws.Range("T9").Value = 0.1144850175826347

# Row 10
$ws.Range("G10").Value = 0.894276
$ws.Range("H10").Value = 2.682828
$ws.Range("I10").Value = 0.1947545739044081
$ws.Range("J10").Value = 0.194754573904408
$ws.Range("M10").Value = 6.492744666666667
$ws.Range("N10").Value = 19.478234
$ws.Range("O10").Value = 0.284509526105254
$ws.Range("P10").Value = 0.2845095261052539
$ws.Range("Q10").Value = 5.806305729528
$ws.Range("R10").Value = 52.256751565752
$ws.Range("S10").Value = 0.0554095315283738
$ws.Range("T10").Value = 0.05540953152837378

# Row 11
$ws.Range("G11").Value = 0.894276
$ws.Range("H11").Value = 2.682828
$ws.Range("I11").Value = 0.1947545739044081
$ws.Range("J11").Value = 0.194754573904408
$ws.Range("O11").Value = 0.1833546924709238
$ws.Range("P11").Value = 0.1833546924709237
$ws.Range("Q11").Value = 3.741925326732
$ws.Range("R11").Value = 33.677327940588
$ws.Range("S11").Value = 0.03570916500554853
$ws.Range("T11").Value = 0.03570916500554852

# Row 12
$ws.Range("G12").Value = 0.894276
$ws.Range("H12").Value = 2.682828
$ws.Range("I12").Value = 0.1947545739044081
$ws.Range("J12").Value = 0.194754573904408
$ws.Range("M12").Value = 1.79534
$ws.Range("N12").Value = 5.38602
$ws.Range("O12").Value = 0.07867109501782452
$ws.Range("P12").Value = 0.0786710950178245
$ws.Range("Q12").Value = 1.60552947384
$ws.Range("R12").Value = 14.44976526456
$ws.Range("S12").Value = 0.01532155558878961
$ws.Range("T12").Value = 0.01532155558878961

# Row 13
$ws.Range("G13").Value = 0.894276
$ws.Range("H13").Value = 2.682828
$ws.Range("I13").Value = 0.1947545739044081
$ws.Range("J13").Value = 0.194754573904408
$ws.Range("M13").Value = 10.34844233333333
$ws.Range("N13").Value = 31.045327
$ws.Range("O13").Value = 0.4534646864059979
$ws.Range("P13").Value = 0.4534646864059979
$ws.Range("Q13").Value = 9.254363616083999
$ws.Range("R13").Value = 83.28927254475599
$ws.Range("S13").Value = 0.08831432178169614
$ws.Range("T13").Value = 0.08831432178169613

# Row 14
$ws.Range("G14").Value = 1.396195
$ws.Range("H14").Value = 4.188585
$ws.Range("I14").Value = 0.3040620147610637
$ws.Range("J14").Value = 0.3040620147610637
$ws.Range("M14").Value = 6.492744666666667
$ws.Range("N14").Value = 19.478234
$ws.Range("O14").Value = 0.284509526105254
$ws.Range("P14").Value = 0.2845095261052539
$ws.Range("Q14").Value = 9.065137639876665
$ws.Range("R14").Value = 81.58623875889
$ws.Range("S14").Value = 0.08650853972627898
$ws.Range("T14").Value = 0.08650853972627896

# Row 15
$ws.Range("G15").Value = 1.396195
$ws.Range("H15").Value = 4.188585
$ws.Range("I15").Value = 0.3040620147610637
$ws.Range("J15").Value = 0.3040620147610637
$ws.Range("O15").Value = 0.1833546924709238
$ws.Range("P15").Value = 0.1833546924709237
$ws.Range("Q15").Value = 5.842108511865
$ws.Range("R15").Value = 52.57897660678501
$ws.Range("S15").Value = 0.05575119720860432
$ws.Range("T15").Value = 0.05575119720860431

# Row 16
$ws.Range("G16").Value = 1.396195
$ws.Range("H16").Value = 4.188585
$ws.Range("I16").Value = 0.3040620147610637
$ws.Range("J16").Value = 0.3040620147610637
$ws.Range("M16").Value = 1.79534
$ws.Range("N16").Value = 5.38602
$ws.Range("O16").Value = 0.07867109501782452
$ws.Range("P16").Value = 0.0786710950178245
$ws.Range("Q16").Value = 2.5066447313
$ws.Range("R16").Value = 22.5598025817
$ws.Range("S16").Value = 0.0239208916545788
$ws.Range("T16").Value = 0.0239208916545788

# Row 17
$ws.Range("G17").Value = 1.396195
$ws.Range("H17").Value = 4.188585
$ws.Range("I17").Value = 0.3040620147610637
$ws.Range("J17").Value = 0.3040620147610637
$ws.Range("M17").Value = 10.34844233333333
$ws.Range("N17").Value = 31.045327
$ws.Range("O17").Value = 0.4534646864059979
$ws.Range("P17").Value = 0.4534646864059979
$ws.Range("Q17").Value = 14.44844344358833
$ws.Range("R17").Value = 130.035990992295
$ws.Range("S17").Value = 0.1378813861716016
$ws.Range("T17").Value = 0.1378813861716016
